$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7663920000000001
$ws.Range("H2").Value = 2.299176
$ws.Range("I2").Value = 0.4782815633346924
$ws.Range("J2").Value = 0.4782815633346925
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 122.226983783168
$ws.Range("R2").Value = 1100.042854048512
$ws.Range("S2").Value = 0.1426850254280123
$ws.Range("T2").Value = 0.1426850254280123
$ws.Range("G3").Value = 0.7663920000000001
$ws.Range("H3").Value = 2.299176
$ws.Range("I3").Value = 0.4782815633346924
$ws.Range("J3").Value = 0.4782815633346925
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 132.247451632824
$ws.Range("R3").Value = 1190.227064695416
$ws.Range("S3").Value = 0.1543826936979351
$ws.Range("T3").Value = 0.1543826936979351
$ws.Range("G4").Value = 0.7663920000000001
$ws.Range("H4").Value = 2.299176
$ws.Range("I4").Value = 0.4782815633346924
$ws.Range("J4").Value = 0.4782815633346925
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 57.010142521288
$ws.Range("R4").Value = 513.091282691592
$ws.Range("S4").Value = 0.06655235516353131
$ws.Range("T4").Value = 0.06655235516353133
$ws.Range("G5").Value = 0.7663920000000001
$ws.Range("H5").Value = 2.299176
$ws.Range("I5").Value = 0.4782815633346924
$ws.Range("J5").Value = 0.4782815633346925
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 44.768493108152
$ws.Range("R5").Value = 402.916437973368
$ws.Range("S5").Value = 0.05226172960990737
$ws.Range("T5").Value = 0.05226172960990738
$ws.Range("G6").Value = 0.7663920000000001
$ws.Range("H6").Value = 2.299176
$ws.Range("I6").Value = 0.4782815633346924
$ws.Range("J6").Value = 0.4782815633346925
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 53.45294197267201
$ws.Range("R6").Value = 481.076477754048
$ws.Range("S6").Value = 0.0623997594353063
$ws.Range("T6").Value = 0.06239975943530632
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8359946666666667
$ws.Range("H7").Value = 2.507984
$ws.Range("I7").Value = 0.5217184366653075
$ws.Range("J7").Value = 0.5217184366653076
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 133.3274702312676
$ws.Range("R7").Value = 1199.947232081408
$ws.Range("S7").Value = 0.1556434830622135
$ws.Range("T7").Value = 0.1556434830622135
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8359946666666667
$ws.Range("H8").Value = 2.507984
$ws.Range("I8").Value = 0.5217184366653075
$ws.Range("J8").Value = 0.5217184366653076
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 144.257983180016
$ws.Range("R8").Value = 1298.321848620144
$ws.Range("S8").Value = 0.1684035174650927
$ws.Range("T8").Value = 0.1684035174650928
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8359946666666667
$ws.Range("H9").Value = 2.507984
$ws.Range("I9").Value = 0.5217184366653075
$ws.Range("J9").Value = 0.5217184366653076
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 62.18772520290311
$ws.Range("R9").Value = 559.689526826128
$ws.Range("S9").Value = 0.07259654846451681
$ws.Range("T9").Value = 0.07259654846451684
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8359946666666667
$ws.Range("H10").Value = 2.507984
$ws.Range("I10").Value = 0.5217184366653075
$ws.Range("J10").Value = 0.5217184366653076
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 48.83430603805689
$ws.Range("R10").Value = 439.508754342512
$ws.Range("S10").Value = 0.05700806796607737
$ws.Range("T10").Value = 0.05700806796607739
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8359946666666667
$ws.Range("H11").Value = 2.507984
$ws.Range("I11").Value = 0.5217184366653075
$ws.Range("J11").Value = 0.5217184366653076
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 58.30746459618133
$ws.Range("R11").Value = 524.767181365632
$ws.Range("S11").Value = 0.06806681970740701
$ws.Range("T11").Value = 0.06806681970740704
